# Roster refresh: rewrite the player/position/team table (18 -> 19 rows).
# Values are written column-by-column (all of A, then all of B, then all of
# C) so that newly-introduced strings are interned in the same relative
# order the source workbook uses.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged, kept for completeness)
$ws.Range("A1").Value = "Oyuncu Adı"
$ws.Range("B1").Value = "Pozisyon"
$ws.Range("C1").Value = "Takım"

# Column A (player names) - new strings introduced here appear in this order
$ws.Range("A2").Value = "Luke Kennard"
$ws.Range("A3").Value = "Deandre Ayton"
$ws.Range("A4").Value = "LaMelo Ball"
$ws.Range("A5").Value = "Devin Vassell"
$ws.Range("A6").Value = "Damian Lillard"
$ws.Range("A7").Value = "Derrick White"
$ws.Range("A8").Value = "Alex Caruso"
$ws.Range("A9").Value = "Onyeka Okongwu"
$ws.Range("A10").Value = "Malik Monk"
$ws.Range("A11").Value = "Keon Ellis"
$ws.Range("A12").Value = "Anthony Davis"
$ws.Range("A13").Value = "Collin Sexton"
$ws.Range("A14").Value = "Cade Cunningham"
$ws.Range("A15").Value = "Julius Randle"
$ws.Range("A16").Value = "Naz Reid"
$ws.Range("A17").Value = "Coby White"
$ws.Range("A18").Value = "Isaiah Hartenstein"
$ws.Range("A19").Value = "Cameron Johnson"

# Column B (positions)
$ws.Range("B2").Value = "SG"
$ws.Range("B3").Value = "C"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("B6").Value = "PG"
$ws.Range("B7").Value = "PG,SG"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("B9").Value = "PF,C"
$ws.Range("B10").Value = "PG,SG,SF"
$ws.Range("B11").Value = "SG,SF"
$ws.Range("B12").Value = "PF,C"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("B15").Value = "PF,C"
$ws.Range("B16").Value = "PF,C"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("B18").Value = "C"
$ws.Range("B19").Value = "SF,PF"

# Column C (teams)
$ws.Range("C2").Value = "Memphis Grizzlies"
$ws.Range("C3").Value = "Portland Trail Blazers"
$ws.Range("C4").Value = "Charlotte Hornets"
$ws.Range("C5").Value = "San Antonio Spurs"
$ws.Range("C6").Value = "Milwaukee Bucks"
$ws.Range("C7").Value = "Boston Celtics"
$ws.Range("C8").Value = "Oklahoma City Thunder"
$ws.Range("C9").Value = "Atlanta Hawks"
$ws.Range("C10").Value = "Sacramento Kings"
$ws.Range("C11").Value = "Sacramento Kings"
$ws.Range("C12").Value = "Los Angeles Lakers"
$ws.Range("C13").Value = "Utah Jazz"
$ws.Range("C14").Value = "Detroit Pistons"
$ws.Range("C15").Value = "Minnesota Timberwolves"
$ws.Range("C16").Value = "Minnesota Timberwolves"
$ws.Range("C17").Value = "Chicago Bulls"
$ws.Range("C18").Value = "Oklahoma City Thunder"
$ws.Range("C19").Value = "Brooklyn Nets"
